# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C for rows 2-43
# from 45795 (2025-05-18) to 45796 (2025-05-19), keeping existing
# number formatting/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C43")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45795) {
        $cell.Value2 = 45796
    }
}
